$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update betting odds / count values for Jogos da Semana FlashScore 2024-10-15
# Each assignment below updates a single cell value to match the new dataset snapshot.

# Row 2
$ws.Range("G2").Value = 1.9
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.63
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("X2").Value = 8
$ws.Range("AH2").Value = 21
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 11
$ws.Range("AS2").Value = 201
$ws.Range("AU2").Value = 9
$ws.Range("AX2").Value = 26
$ws.Range("BA2").Value = 126

# Row 3
$ws.Range("G3").Value = 2.15
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("Y3").Value = 9
$ws.Range("AA3").Value = 17
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 51
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 12
$ws.Range("AK3").Value = 26

# Row 5
$ws.Range("H5").Value = 4.05
$ws.Range("I5").Value = 1.39
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 2.22
$ws.Range("L5").Value = 1.93
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 7.5
$ws.Range("O5").Value = 1.28
$ws.Range("P5").Value = 3.35
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 1.88
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.72
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.65
$ws.Range("W5").Value = 18
$ws.Range("AA5").Value = 100
$ws.Range("AB5").Value = 90
$ws.Range("AC5").Value = 7.5
$ws.Range("AD5").Value = 8.25
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 110
$ws.Range("AG5").Value = 6
$ws.Range("AH5").Value = 6
$ws.Range("AK5").Value = 12
$ws.Range("AL5").Value = 32
$ws.Range("AP5").Value = 50
$ws.Range("AR5").Value = 400
$ws.Range("AT5").Value = 2.72
$ws.Range("AU5").Value = 8.75
$ws.Range("AV5").Value = 90
$ws.Range("AW5").Value = 3.1
$ws.Range("AX5").Value = 6.6
$ws.Range("AY5").Value = 18.5
$ws.Range("AZ5").Value = 20
$ws.Range("BA5").Value = 55
$ws.Range("BB5").Value = 300

# Row 6
$ws.Range("G6").Value = 1.31
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 7.9
$ws.Range("K6").Value = 2.52
$ws.Range("L6").Value = 7
$ws.Range("V6").Value = 1.9
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 7.2
$ws.Range("Z6").Value = 8.75
$ws.Range("AG6").Value = 25
$ws.Range("AH6").Value = 60
$ws.Range("AI6").Value = 25
$ws.Range("AK6").Value = 90
$ws.Range("AL6").Value = 70
$ws.Range("AO6").Value = 5.8
$ws.Range("AW6").Value = 9.25
$ws.Range("BB6").Value = 450

# Row 7
$ws.Range("G7").Value = 5.1
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 1.57
$ws.Range("J7").Value = 5.3
$ws.Range("K7").Value = 2.22
$ws.Range("L7").Value = 2.12
$ws.Range("N7").Value = 7.5
$ws.Range("O7").Value = 1.28
$ws.Range("P7").Value = 3.35
$ws.Range("Q7").Value = 1.85
$ws.Range("R7").Value = 1.9
$ws.Range("S7").Value = 1.38
$ws.Range("T7").Value = 2.82
$ws.Range("W7").Value = 13.5
$ws.Range("X7").Value = 30
$ws.Range("Y7").Value = 17
$ws.Range("Z7").Value = 100
$ws.Range("AA7").Value = 55
$ws.Range("AB7").Value = 55
$ws.Range("AC7").Value = 7.5
$ws.Range("AD7").Value = 7.4
$ws.Range("AE7").Value = 17
$ws.Range("AG7").Value = 6.7
$ws.Range("AH7").Value = 7.3
$ws.Range("AI7").Value = 8
$ws.Range("AJ7").Value = 11.5
$ws.Range("AK7").Value = 12.5
$ws.Range("AL7").Value = 27
$ws.Range("AN7").Value = 6.8
$ws.Range("AO7").Value = 30
$ws.Range("AP7").Value = 35
$ws.Range("AQ7").Value = 200
$ws.Range("AR7").Value = 250
$ws.Range("AS7").Value = 500
$ws.Range("AT7").Value = 2.82
$ws.Range("AU7").Value = 7.8
$ws.Range("AW7").Value = 3.4
$ws.Range("AX7").Value = 7.6
$ws.Range("AY7").Value = 17.5
$ws.Range("AZ7").Value = 24
$ws.Range("BA7").Value = 55
